$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$sh.TextFrame.TextRange.Text = "NEWTestest12"
